$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "43.047.05"
$ws.Range("E2").Value = "  +0.54%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.336.19"
$ws.Range("E3").Value = "  +4.29%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.17%  "

# Row 5: Solana -> BNB
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'310.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.35%  "

# Row 6: BNB -> Solana
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'108.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.01%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.06%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.40%  "

# Row 9: Cardano
$ws.Range("D9").Value = "'0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "

# Row 10: Avalanche
$ws.Range("D10").Value = "'43.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.40%  "

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.0933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "

# Row 12: Polkadot
$ws.Range("D12").Value = "'8.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "

# Row 13: Polygon
$ws.Range("D13").Value = "'1.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +22.12%  "

# Row 14: TRON
$ws.Range("E14").Value = "  -0.13%  "

# Row 15: Chainlink
$ws.Range("D15").Value = "'16.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.81%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.681.05"
$ws.Range("E16").Value = "  +3.92%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.387.80"
$ws.Range("E17").Value = "  +5.75%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "43.057.78"
$ws.Range("E18").Value = "  +0.73%  "

# Row 19: ShibaInu
$ws.Range("E19").Value = "  +0.29%  "

# Row 20: Uniswap
$ws.Range("D20").Value = "'7.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.75%  "

# Row 21: Litecoin
$ws.Range("D21").Value = "'75.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.70%  "

# Row 22: PancakeSwap
$ws.Range("D22").Value = "'3.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.24%  "

# Row 23: ImmutableX
$ws.Range("D23").Value = "'2.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.80%  "

# Row 24: BitcoinCash
$ws.Range("D24").Value = "'252.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.15%  "

# Row 25: InternetComputer(DFINITY)
$ws.Range("D25").Value = "'9.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.58%  "

# Row 26: Cosmos
$ws.Range("D26").Value = "'11.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.18%  "

# Row 27: Dai
$ws.Range("E27").Value = "  -0.03%  "

# Row 28: InjectiveProtocol
$ws.Range("D28").Value = "'39.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.16%  "

# Row 29: Toncoin
$ws.Range("E29").Value = "  +0.87%  "

# Row 30: EthereumClassic
$ws.Range("D30").Value = "'22.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.35%  "

# Row 31: Monero
$ws.Range("D31").Value = "'173.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.46%  "

# Row 32: WEMIXToken
$ws.Range("D32").Value = "'3.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.36%  "

# Row 33: Hedera
$ws.Range("D33").Value = "'0.0910"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.53%  "

# Row 34: Filecoin
$ws.Range("D34").Value = "'5.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.49%  "

# Row 35: RenderToken
$ws.Range("D35").Value = "'5.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "

# Row 36: Stellar
$ws.Range("E36").Value = "  +2.45%  "

# Row 37: NEARProtocol -> VeChain
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0378"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.96%  "

# Row 38: VeChain -> NEARProtocol
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.80%  "

# Row 39: Kaspa
$ws.Range("E39").Value = "  -0.98%  "

# Row 40: LidoDAOToken
$ws.Range("E40").Value = "  +7.18%  "

# Row 41: ARBITRUM
$ws.Range("E41").Value = "  +12.59%  "

# Row 42: MultiversX
$ws.Range("D42").Value = "'71.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.27%  "

# Row 43: Algorand
$ws.Range("E43").Value = "  -1.74%  "

# Row 44: FirstDigitalUSD
$ws.Range("E44").Value = "  -0.20%  "

# Row 45: Celestia
$ws.Range("E45").Value = "  -5.69%  "

# Row 46: THORChain
$ws.Range("D46").Value = "'5.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.42%  "

# Row 47: Aave -> FraxShare
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.26%  "

# Row 48: FraxShare -> Aave
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'110.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.30%  "

# Row 49: TrustWalletToken
$ws.Range("D49").Value = "'1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.32%  "

# Row 50: Cronos
$ws.Range("D50").Value = "'0.0996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "

# Row 51: ordi
$ws.Range("D51").Value = "'70.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.71%  "
